# timeLog.xlsx maintenance edit:
#  - bring the "second" time table (columns J:N) in line with the formatting
#    already used by the first table (columns A:E) by painting the same
#    per-column cell formats down over rows 12-24
#  - log a missed day (2024-04-02) that had already happened by the time of
#    the previous entry (row 14, J:N table) and also add it as a fresh row
#    at the bottom of the first table (row 32)
#  - likewise normalize the formatting of the already-entered late rows
#    (26-31) of the first table
#  - move the active selection to where the user was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- 1. Normalize formatting of the J:N ("Viki") table rows 12-24 to match
#        the already-clean formatting used in rows 3-11 ---------------------
$ws.Range("J3:N3").Copy() | Out-Null
$ws.Range("J12:N24").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Fill in the forgotten entry for 2024-04-02 in the J:N table (row 14)
$ws.Range("J14").Value = "04/02/2024"
$ws.Range("K14").Formula = "=8"
$ws.Range("L14").Formula = "=9+35/60"

# --- 3. Normalize formatting of the late A:E rows (26-31) to match rows 3-11
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A26:E31").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- 4. Append the same 2024-04-02 entry to the main (A:E) table as row 32 -
$ws.Range("A26:E26").Copy() | Out-Null
$ws.Range("A32:E32").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A32").Value = "04/02/2024"
$ws.Range("B32").Formula = "=8"
$ws.Range("C32").Formula = "=9+35/60"
$ws.Range("D32").ClearContents() | Out-Null
$ws.Range("E32").ClearContents() | Out-Null

# --- 5. Leave the selection where the user ended up working -----------------
$ws.Range("G17").Select() | Out-Null
